# Fix duplicate/misclassified entry for "Bjarne Forfot".
# He was incorrectly classified as Female; move his two result rows
# (25m and 50m pools) from the Female sheets into the correct position
# on the Male sheets (sorted by Poeng/points, descending).

$wb = $excel.ActiveWorkbook

# --- Male_25m: insert Bjarne Forfot's 25m result at row 122 ---
$wsMale25 = $wb.Worksheets.Item("Male_25m")
$wsMale25.Rows.Item(122).Insert()
$wsMale25.Cells.Item(122, 1).Value = "Bjarne Forfot"
$wsMale25.Cells.Item(122, 2).Value = "1.24,28"
$wsMale25.Cells.Item(122, 3).Value = 188
$wsMale25.Cells.Item(122, 4).Value = "30.09.2017"
$wsMale25.Cells.Item(122, 5).Value = "Husebybadet"
$wsMale25.Cells.Item(122, 6).Value = "25m"
$wsMale25.Cells.Item(122, 7).Value = "Male"

# --- Male_50m: insert Bjarne Forfot's 50m result at row 76 ---
$wsMale50 = $wb.Worksheets.Item("Male_50m")
$wsMale50.Rows.Item(76).Insert()
$wsMale50.Cells.Item(76, 1).Value = "Bjarne Forfot"
$wsMale50.Cells.Item(76, 2).Value = "1.30,29"
$wsMale50.Cells.Item(76, 3).Value = 186
$wsMale50.Cells.Item(76, 4).Value = "14.04.2018"
$wsMale50.Cells.Item(76, 5).Value = "Bergen"
$wsMale50.Cells.Item(76, 6).Value = "50m"
$wsMale50.Cells.Item(76, 7).Value = "Male"

# --- Female_25m: remove the erroneous Bjarne Forfot row (row 136) ---
$wsFemale25 = $wb.Worksheets.Item("Female_25m")
$wsFemale25.Rows.Item(136).Delete()

# --- Female_50m: remove the erroneous Bjarne Forfot row (row 88) ---
$wsFemale50 = $wb.Worksheets.Item("Female_50m")
$wsFemale50.Rows.Item(88).Delete()
